$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "60.186.51"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +2.70%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.634.04"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +0.16%  "

$ws.Range("E4").Value = "  +0.09%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "567.96"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +5.90%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "145.61"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.90%  "

$ws.Range("E7").Value = "  -0.22%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "6.85"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -2.42%  "

$ws.Range("E10").Value = "  +3.62%  "

$ws.Range("E11").Value = "  +6.24%  "

$ws.Range("E12").Value = "  +2.29%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "3.105.16"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.13%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "60.182.43"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +2.79%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "21.73"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +4.04%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "2.645.14"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +1.00%  "

$ws.Range("E17").Value = "  +2.85%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "4.58"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +3.69%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "343.47"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +2.64%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "10.39"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.13%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.34"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.36%  "

$ws.Range("E22").Value = "  +0.11%  "

$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("E24").Value = "  +5.15%  "

$ws.Range("E25").Value = "  +1.80%  "

$ws.Range("E26").Value = "  -0.20%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "7.31"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +2.12%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.0₃0772"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +4.26%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.09%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.71"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +3.93%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "6.12"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +5.07%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "156.31"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +3.79%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "19.18"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.27%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.09"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +4.92%  "

$ws.Range("E35").Value = "  +8.23%  "

$ws.Range("E36").Value = "  +12.35%  "

$ws.Range("E37").Value = "  +5.31%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "37.50"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.90%  "

$ws.Range("E39").Value = "  +5.56%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "303.10"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +7.53%  "

$ws.Range("E41").Value = "  +2.27%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.994"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.43%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.605"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.70%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.0974"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +4.14%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0548"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +3.11%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "19.32"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +1.11%  "

$ws.Range("E47").Value = "  -0.51%  "

$ws.Range("E48").Value = "  +4.91%  "

$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "123.37"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +10.06%  "

$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.959.42"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.58%  "

$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "4.59"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +2.97%  "
